# Generate Report for Handback
# Updates the localization-status workbook to reflect that the two
# in-flight files (de8179a3-c92e-43be-a73e-685c667cbc87.md) have now been
# handed back: the "Status" column moves from "Ready for handoff" to
# "Handed back: in sync with en-US", and the "Latest Handback DateTime"
# column is stamped with the new handback time for each locale sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to de8179a3-...md, columns B (zh-cn) and C (de-de)
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is the de8179a3-...md entry
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G3").Value = "2016-03-09 04:59:45"

# de-de sheet: row 3 is the de8179a3-...md entry
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G3").Value = "2016-03-09 04:59:50"
